$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.096435308456421
$ws.Range("B1").Value = 1.899765133857727
$ws.Range("D1").Value = 1.107171177864075
$ws.Range("E1").Value = 1.135149240493774
